# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells that get refreshed each time
# the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# "Overview" sheet: Latest HO Xliff Generate Date (shared with de-de!H2 below)
$wsOverview.Range("G2").Value = "2016-09-06 23:21:12"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-06 23:20:59"
$wsZhCn.Range("K2").Value = "2016-09-06 23:21:33"

# "de-de" sheet: Correspond Handoff Datetime (mirrors Overview!G2) /
# Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-06 23:21:12"
$wsDeDe.Range("K2").Value = "2016-09-06 23:21:41"
